# P4-2136 test data obfuscating exercise.
# Renames the "Bedford" test-data location (County Court + Prison records)
# to a fictitious "Fred" location across the Courts and Prisons sheets.

$wb = $excel.ActiveWorkbook

# ---- Courts sheet: "Bedford County Court" row (row 2) ----
$courts = $wb.Worksheets.Item("Courts")
$courts.Range("C2").Value = "Freds County Court"
$courts.Range("D2").Value = "FRDCT"
$courts.Range("F2").Value = "Freds County Court"
$courts.Range("J2").Value = "Fredford"
$courts.Range("K2").Value = "Fredfordshire"
$courts.Range("L2").Value = "England"
$courts.Range("M2").Value = "FR31 3ZZ"

# ---- Prisons sheet: "HMP Bedford" row (row 2) ----
$prisons = $wb.Worksheets.Item("Prisons")
$prisons.Range("C2").Value = "HMP Fred"
$prisons.Range("F2").Value = "HMP Fred"
$prisons.Range("H2").Value = "Fred"
$prisons.Range("I2").Value = "Fredfordshire"
$prisons.Range("J2").Value = "Fred St"
$prisons.Range("K2").Value = "FD40 1HG"

# ---- View state: bring the Prisons sheet to the front, selected at K2 ----
$courts.Range("C2").Select() | Out-Null

$prisons.Activate() | Out-Null
$prisons.Range("K2").Select() | Out-Null
